# 10Th - MB for single stock and added new group
#
# The sheet is a rolling-window "MarketBeat rank" report: column A holds
# analyst/firm names, and each of the other columns is one date snapshot
# (header = "Jun_NN") with "UN" (unchanged) in every row except where a
# rating-change note is recorded.
#
# This edit rolls the date window forward:
#   - drops the two oldest "UN"-only snapshot columns (Jun_24, Jun_19)
#   - keeps the remaining snapshot columns, shifting left
#   - adds two new, newest snapshot columns (Jun_27, Jun_26) at the front
#   - appends two new analyst/firm rows (Benchmark, Evercore ISI)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Columns -----------------------------------------------------------
# Before: A=names B=Jun_24(UN) C=Jun_22 D=Jun_19(UN) E=Jun_17 F=Jun_15 G=Jun_13 H=Jun_10
# Drop the two purely-"UN" snapshot columns (no rating notes live there).
$ws.Columns("D").Delete() | Out-Null   # remove Jun_19
$ws.Columns("B").Delete() | Out-Null   # remove Jun_24
# Now: A=names B=Jun_22(notes) C=Jun_17 D=Jun_15(notes) E=Jun_13 F=Jun_10(notes)

# The old "Jun_22" column becomes "Jun_26" (its data/notes stay put).
$ws.Range("B1").Value = "Jun_26"

# Insert a brand-new newest-snapshot column in front of it.
$ws.Columns("B").Insert() | Out-Null
# Now: A=names B=(new) C=Jun_26(notes) D=Jun_17 E=Jun_15(notes) F=Jun_13 G=Jun_10(notes)

$ws.Range("B1").Value = "Jun_27"
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
}

# --- New rows ------------------------------------------------------------
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
